$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right 5 -> 4, Wrong -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total): Right 40 -> 32, Wrong -10 -> -20, Max text "40 / 140" -> "12 / 112"
$ws.Range("B12").Value = 32
$ws.Range("C12").Value = -20
$ws.Range("E12").Value = "12 / 112"
